# Update backup workbook:
#  1. Seat Assignments: remove the seat-assignment row for the contestant who
#     has become a standby (row r=2 in the original sheet, uuid a001b03c...),
#     shifting the remaining rows up. The row that ends up at r=2 also picks
#     up an (empty) Notes cell in H2, matching the source data export.
#  2. Insert a new "Standbys" worksheet right after "Seat Assignments" (and
#     before "Groups") holding the standby record that used to be seated.

$wb = $excel.ActiveWorkbook

# --- 1. Seat Assignments: delete the now-obsolete seat row -----------------
$seatAssignments = $wb.Worksheets.Item("Seat Assignments")
$seatAssignments.Rows.Item(2).Delete()

# Recreate the blank (but present/typed) Notes cell at H2 the export produces
# for this row -- a real empty-text cell, not just an absent one.
$h2 = $seatAssignments.Cells.Item(2, 8)
$h2.Formula = "'"
$h2.ClearFormats()

# --- 2. New "Standbys" sheet, inserted between "Seat Assignments" and
#        "Groups" -------------------------------------------------------
$standbys = $wb.Worksheets.Add($null, $seatAssignments)
$standbys.Name = "Standbys"

$standbys.Range("A1").Value = "ID"
$standbys.Range("B1").Value = "RecordDayID"
$standbys.Range("C1").Value = "ContestantID"
$standbys.Range("D1").Value = "Status"
$standbys.Range("E1").Value = "Notes"

$standbys.Range("A2").Value = "cc98ff1c-1a1d-43f8-82ea-abcca3482b1e"
$standbys.Range("B2").Value = "e432f0fe-1383-44a2-990c-5f787da5008a"
$standbys.Range("C2").Value = "d698b1de-6641-45c6-aa63-f577d2b634bb"
$standbys.Range("D2").Value = "pending"
